# Fruta / hortaliza, semanal
# Insert one new weekly record for "Macroferia Regional de Talca - Berenjena"
# at row 87 (pushing the existing rows 87-115 down to 88-116).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 87; this shifts rows
# 87..115 down to 88..116 and grows the sheet dimension to A1:R116.
$ws.Rows(87).Insert()

# Populate the newly inserted row 87 with the new weekly data point.
$ws.Cells.Item(87, 1).Value = 5
$ws.Cells.Item(87, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(87, 3).Value = "Maule"
$ws.Cells.Item(87, 4).Value = 44627
$ws.Cells.Item(87, 5).Value = 7
$ws.Cells.Item(87, 6).Value = 100112001
$ws.Cells.Item(87, 7).Value = "Berenjena"
$ws.Cells.Item(87, 8).Value = "Sin especificar"
$ws.Cells.Item(87, 9).Value = "Primera"
$ws.Cells.Item(87, 10).Value = 150
$ws.Cells.Item(87, 11).Value = 7000
$ws.Cells.Item(87, 12).Value = 7000
$ws.Cells.Item(87, 13).Value = 7000
$ws.Cells.Item(87, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(87, 15).Value = "Región del Maule"
$ws.Cells.Item(87, 16).Value = 140
$ws.Cells.Item(87, 17).Value = 50
$ws.Cells.Item(87, 18).Value = "Hortaliza"
